# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns with
# the latest snapshot values for each coin row on Sheet1.
#
# A handful of "Price" values (column D) are plain decimal-looking strings
# (e.g. "235.34", "0.388", "58.75") that Excel's COM layer would otherwise
# auto-convert to numbers on assignment (losing formatting like trailing
# zeros, e.g. "0.130" -> 0.13). Those cells are written with a leading
# apostrophe to force text, then immediately restyled back to "Normal" so
# no stray number-format/quote-prefix styling is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "37.631.61"; ForceText = $false }
    @{ Cell = "E2"; Value = "  +2.50%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "2.083.03"; ForceText = $false }
    @{ Cell = "E3"; Value = "  +3.71%  "; ForceText = $false }
    @{ Cell = "E4"; Value = "  -0.08%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "235.34"; ForceText = $true }
    @{ Cell = "E5"; Value = "  -0.41%  "; ForceText = $false }
    @{ Cell = "E6"; Value = "  +4.35%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "58.75"; ForceText = $true }
    @{ Cell = "E7"; Value = "  +6.33%  "; ForceText = $false }
    @{ Cell = "E8"; Value = "  -0.03%  "; ForceText = $false }
    @{ Cell = "D9"; Value = "0.388"; ForceText = $true }
    @{ Cell = "E9"; Value = "  +4.16%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "58.81"; ForceText = $true }
    @{ Cell = "E10"; Value = "  +1.16%  "; ForceText = $false }
    @{ Cell = "D11"; Value = "0.0764"; ForceText = $true }
    @{ Cell = "E11"; Value = "  +2.18%  "; ForceText = $false }
    @{ Cell = "E12"; Value = "  +3.50%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "2.385.86"; ForceText = $false }
    @{ Cell = "E13"; Value = "  +3.56%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "14.59"; ForceText = $true }
    @{ Cell = "E14"; Value = "  +2.35%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "21.22"; ForceText = $true }
    @{ Cell = "E15"; Value = "  +4.74%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "0.783"; ForceText = $true }
    @{ Cell = "E16"; Value = "  +3.23%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "5.23"; ForceText = $true }
    @{ Cell = "E17"; Value = "  +2.10%  "; ForceText = $false }
    @{ Cell = "D18"; Value = "2.079.47"; ForceText = $false }
    @{ Cell = "E18"; Value = "  +3.16%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "37.750.14"; ForceText = $false }
    @{ Cell = "E19"; Value = "  +3.08%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "6.23"; ForceText = $true }
    @{ Cell = "E20"; Value = "  +16.55%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "70.02"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +3.17%  "; ForceText = $false }
    @{ Cell = "E22"; Value = "  +1.35%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "227.11"; ForceText = $true }
    @{ Cell = "E23"; Value = "  +2.24%  "; ForceText = $false }
    @{ Cell = "E24"; Value = "  -0.03%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "2.50"; ForceText = $true }
    @{ Cell = "E25"; Value = "  +3.23%  "; ForceText = $false }
    @{ Cell = "E26"; Value = "  +0.60%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "167.56"; ForceText = $true }
    @{ Cell = "E27"; Value = "  +2.81%  "; ForceText = $false }
    @{ Cell = "E28"; Value = "  +9.32%  "; ForceText = $false }
    @{ Cell = "D29"; Value = "9.07"; ForceText = $true }
    @{ Cell = "E29"; Value = "  +4.38%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "0.130"; ForceText = $true }
    @{ Cell = "E30"; Value = "  +1.08%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "19.31"; ForceText = $true }
    @{ Cell = "E31"; Value = "  +2.14%  "; ForceText = $false }
    @{ Cell = "E32"; Value = "  +1.88%  "; ForceText = $false }
    @{ Cell = "D33"; Value = "4.55"; ForceText = $true }
    @{ Cell = "E33"; Value = "  +3.93%  "; ForceText = $false }
    @{ Cell = "D34"; Value = "0.0628"; ForceText = $true }
    @{ Cell = "E34"; Value = "  +3.58%  "; ForceText = $false }
    @{ Cell = "E35"; Value = "  +6.79%  "; ForceText = $false }
    @{ Cell = "E36"; Value = "  +8.06%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "3.36"; ForceText = $true }
    @{ Cell = "E38"; Value = "  -0.46%  "; ForceText = $false }
    @{ Cell = "D39"; Value = "5.96"; ForceText = $true }
    @{ Cell = "E39"; Value = "  +3.34%  "; ForceText = $false }
    @{ Cell = "E40"; Value = "  -0.12%  "; ForceText = $false }
    @{ Cell = "D41"; Value = "4.61"; ForceText = $true }
    @{ Cell = "E41"; Value = "  +20.46%  "; ForceText = $false }
    @{ Cell = "E42"; Value = "  -0.92%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "0.0963"; ForceText = $true }
    @{ Cell = "E43"; Value = "  +3.90%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "1.474.14"; ForceText = $false }
    @{ Cell = "E44"; Value = "  +0.88%  "; ForceText = $false }
    @{ Cell = "D45"; Value = "1.19"; ForceText = $true }
    @{ Cell = "E45"; Value = "  +7.21%  "; ForceText = $false }
    @{ Cell = "D46"; Value = "96.20"; ForceText = $true }
    @{ Cell = "E46"; Value = "  +6.40%  "; ForceText = $false }
    @{ Cell = "E47"; Value = "  +4.66%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "15.92"; ForceText = $true }
    @{ Cell = "E48"; Value = "  +3.35%  "; ForceText = $false }
    @{ Cell = "E49"; Value = "  +4.09%  "; ForceText = $false }
    @{ Cell = "E50"; Value = "  +5.71%  "; ForceText = $false }
    @{ Cell = "E51"; Value = "  +1.70%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $range.Value = "'" + $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
